$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "name" column (C) needed to widen to fit its longest translated value
# ("Animal de Estimação"), splitting it out of the former A:D shared-width
# column group.
$ws.Columns.Item(3).ColumnWidth = 19

# Reviewer's cursor ended up on E30 after the pass over the translations.
$ws.Range("E30").Select() | Out-Null
